# Apply the "Continuing cleaning and adding in LOOPR results" edit:
#  - Sheet1: clear D5 (drop the stray https://osf.io/dmf62/ text but keep its
#            style) and record ManyLabs2's n_studies = 28 in H4
#  - Add a new "Codebook" worksheet (after Sheet1) describing the variables
#    used in the replication-results data
#  - Leave Sheet1 as the active sheet with C7 selected; Codebook's own
#    selection lands on E16

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1 data edits -----------------------------------------------------

# D5 previously held "https://osf.io/dmf62/"; clear the text but keep formatting.
$ws1.Range("D5").Value = ""

# ManyLabs2 (row 4) gains an n_studies count.
$ws1.Range("H4").Value = 28

# --- New Codebook sheet -----------------------------------------------------

$codebook = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$codebook.Name = "Codebook"

# Column A (variable names), rows 2-19 first ...
$varNames = @(
    "authorsTitle.o",
    "correlation.o",
    "fis.o",
    "seFish.o",
    "n.o",
    "pVal.o",
    "resultUsedInRep.o",
    "correlation.r",
    "fis.r",
    "n.r",
    "seFish.r",
    "pVal.r",
    "seDifference.ro",
    "source",
    "cohenD.o",
    "seCohenD.o",
    "cohenD.r",
    "seCohenD.r"
)
for ($i = 0; $i -lt $varNames.Length; $i++) {
    $codebook.Cells.Item($i + 2, 1).Value = $varNames[$i]
}

# ... then the header row ...
$codebook.Cells.Item(1, 1).Value = "Variable"
$codebook.Cells.Item(1, 2).Value = "Explanation "

# ... then column B (explanations), rows 2-19.
$explanations = @(
    "Original author title (or unique study identifier from large scale replication project)",
    "Original correlation",
    "Original fisher transformed correlation",
    "Original SE fisher transformed correlation",
    "Original sample size",
    "Original p value (as reported in paper)",
    'Original result (e.g., "t(123) = 123, p < .001"',
    "Replicaiton Correlation",
    "Replication fisher transformed correlation",
    "Replication SE fisher transformed correlation",
    "Replication sample size",
    "Replication p value (as reported in paper)",
    'Replication result (e.g., "t(123) = 123, p < .001"',
    "Which replication project this was from",
    "Original effect in Cohen's d",
    "Original study standard error of d",
    "Replication study effect in Cohen's d",
    "Replication study standard error of d"
)
for ($i = 0; $i -lt $explanations.Length; $i++) {
    $codebook.Cells.Item($i + 2, 2).Value = $explanations[$i]
}

# --- Selections / active sheet ---------------------------------------------

$codebook.Range("E16").Select()

$ws1.Activate()
$ws1.Range("C7").Select()
